{"js": "// Apply the MDL letter field updates:\n//  - Letter date: February 09, 2026 -> February 10, 2026\n//  - Street address: 122 Young Street -> 1130 Bluffs Parkway\n//  - City/state/zip: Henderson, NC 27536 -> Canton, GA 30114\n//  - Audit period end date (appears twice): December 31, 2023 -> June 30, 2022\n\nconst body = context.document.body;\n\nasync function replaceAll(searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nawait replaceAll(\"February 09, 2026\", \"February 10, 2026\");\nawait replaceAll(\"122 Young Street\", \"1130 Bluffs Parkway\");\nawait replaceAll(\"Henderson, NC 27536\", \"Canton, GA 30114\");\n// Replaces both the Subject-line mention and the standalone bold date run,\n// since both share the identical literal text and target replacement.\nawait replaceAll(\"December 31, 2023\", \"June 30, 2022\");\n", "ps1": "# Apply the MDL letter field updates:\n#  - Letter date: February 09, 2026 -> February 10, 2026\n#  - Street address: 122 Young Street -> 1130 Bluffs Parkway\n#  - City/state/zip: Henderson, NC 27536 -> Canton, GA 30114\n#  - Audit period end date (appears twice): December 31, 2023 -> June 30, 2022\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $range = $d.Content\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-AllText \"February 09, 2026\" \"February 10, 2026\"\nReplace-AllText \"122 Young Street\" \"1130 Bluffs Parkway\"\nReplace-AllText \"Henderson, NC 27536\" \"Canton, GA 30114\"\n# Replaces both the Subject-line mention and the standalone bold date run,\n# since both share the identical literal text and target replacement.\nReplace-AllText \"December 31, 2023\" \"June 30, 2022\"\n"}
